$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the "Series Code" column (column D) entirely, shifting all
# subsequent columns (E:BB) one position to the left.
$ws.Columns("D").Delete()

# Add a stray single-space label left behind in the now-shifted layout
# (was previously an empty formatted cell in column D).
$ws.Range("C14").Value = " "

# Restore the active selection to what it was left at after the edit.
$ws.Range("C19").Select()
